$wb = $excel.ActiveWorkbook

# Sheet names that contain the event data table that needs updating.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 290
    $ws.Range("F3").Value = 11
    $ws.Range("F5").Value = 910
    $ws.Range("F6").Value = 215
}
